$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New monthly data rows to append below the existing data (rows 46-48)
$data = @(
    @(45901, 0.231, 0.685, 0.098, 0.276, 1.385),
    @(45931, 0.099, 0.581, 0.263, 0.265, 1.568),
    @(45962, 0.106, 0.6,   0.321, 0.247, 1.869)
)

$startRow = 46
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $vals = $data[$i]

    # Copy the date-column formatting (style) from the row above so the
    # new date cell matches the existing date column styling.
    $ws.Range("A" + ($row - 1)).Copy()
    $ws.Range("A" + $row).PasteSpecial(-4122)

    $ws.Cells.Item($row, 1).Value = $vals[0]
    $ws.Cells.Item($row, 2).Value = $vals[1]
    $ws.Cells.Item($row, 3).Value = $vals[2]
    $ws.Cells.Item($row, 4).Value = $vals[3]
    $ws.Cells.Item($row, 5).Value = $vals[4]
    $ws.Cells.Item($row, 6).Value = $vals[5]
}
